$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Add a new "T_KEYWORD" (关键字表) table block to the bottom of the table
# catalogue, mirroring the layout already used by the other table blocks
# (T_USER, T_VISITOR, T_ARTICLE, T_COMMENT).
#
# Rows 35-42 (the T_COMMENT block) are the closest structural template for
# the new block (rows 45-49), so the formatting is copied from there and the
# cell values are then overwritten.
# ---------------------------------------------------------------------------

# Copy formatting for the new rows from the most similar existing block.
$ws.Range("A35:F36").Copy() | Out-Null
$ws.Range("A45:F46").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A37:F39").Copy() | Out-Null
$ws.Range("A47:F49").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Table header row (physical table name / logical name)
$ws.Range("A45").Value = "表物理名称"
$ws.Range("B45").Value = ""
$ws.Range("C45").Value = "T_KEYWORD"
$ws.Range("D45").Value = "逻辑名称"
$ws.Range("E45").Value = "关键字表"
$ws.Range("F45").Value = ""

# Column header row
$ws.Range("A46").Value = "序号"
$ws.Range("B46").Value = "字段名称"
$ws.Range("C46").Value = "数据类型"
$ws.Range("D46").Value = "描述"
$ws.Range("E46").Value = "备注说明"
$ws.Range("F46").Value = "是否索引"

# Field row 1: ID
$ws.Range("A47").Value = 1
$ws.Range("B47").Value = "ID"
$ws.Range("C47").Value = "INT"
$ws.Range("D47").Value = "ID"
$ws.Range("E47").Value = "主键"
$ws.Range("F47").Value = "Y"

# Field row 2: ARTICLEID
$ws.Range("A48").Value = 2
$ws.Range("B48").Value = "ARTICLEID"
$ws.Range("C48").Value = "INT"
$ws.Range("D48").Value = "文档ID"
$ws.Range("E48").Value = "文档ID"
$ws.Range("F48").Value = ""

# Field row 3: KEYWORD
$ws.Range("A49").Value = 3
$ws.Range("B49").Value = "KEYWORD"
$ws.Range("C49").Value = "VARCHAR2（100）"
$ws.Range("D49").Value = "关键字"
$ws.Range("E49").Value = ""
$ws.Range("F49").Value = ""

# Merge the table-name header cell, same as for the other blocks.
$ws.Range("A45:B45").Merge() | Out-Null

# Update the view: scroll so row 22 is at the top and select B47:D49, matching
# the author's final cursor position after typing in the new table.
$ws.Range("B47:D49").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
